$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.087132930755615
$ws.Range("B1").Value = 2.410617351531982
$ws.Range("C1").Value = 2.433518886566162
$ws.Range("D1").Value = 2.912941217422485
$ws.Range("E1").Value = 0.7953464984893799
